$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose full records (columns B..AD) swapped places while the
# id/Div/Date columns (A, C, D) stay anchored to the row position.
$rowPairs = @(
    ,@(215, 216)
    ,@(263, 265)
    ,@(271, 272)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($v1 -ne $v2) {
            $cell1.Value2 = $v2
            $cell2.Value2 = $v1
        }
    }
}
